# CSCI446_Project2_WumpusWorld - AgentRunStatistics.xlsx
# "fire arrow and kill wumpus working (stat tracking as well)"
#
# This updates the "Knowledge Based Agent" sheet's first results table
# (Map Size: 5x5, rows 3-22) with a fresh batch of simulation run
# statistics, and moves the sheet's active cell selection from K14 to D14.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Knowledge Based Agent")

# New run statistics: columns are
#   B = # Decisions Made, C = # Wumpi Killed, D = # Times Fell Into Pit,
#   E = Gold Found? (1=yes), F = # Times Killed by Wumpus,
#   G = # Cells Explored, H = Score
$newStats = @(
    @(7,  0, 0, 0, 0, 6,  -9),
    @(12, 0, 0, 0, 0, 11, -16),
    @(17, 1, 0, 0, 0, 17, 979),
    @(8,  1, 0, 0, 0, 8,  991),
    @(10, 0, 0, 0, 0, 9,  -16),
    @(16, 1, 0, 0, 0, 16, 979),
    @(1,  0, 0, 0, 0, 0,  0),
    @(9,  1, 0, 0, 0, 9,  990),
    @(15, 0, 0, 0, 0, 14, -19),
    @(9,  1, 0, 0, 0, 9,  987),
    @(6,  0, 0, 0, 0, 5,  -8),
    @(12, 1, 0, 0, 0, 12, 984),
    @(3,  0, 0, 0, 0, 2,  -3),
    @(12, 0, 0, 0, 0, 11, -18),
    @(11, 0, 0, 0, 0, 10, -13),
    @(11, 1, 0, 0, 0, 11, 987),
    @(4,  1, 0, 0, 0, 4,  996),
    @(10, 1, 0, 0, 0, 10, 989),
    @(22, 1, 1, 0, 0, 22, 980),
    @(14, 0, 0, 0, 0, 13, -20)
)

$firstRow = 3
for ($i = 0; $i -lt $newStats.Count; $i++) {
    $row  = $firstRow + $i
    $vals = $newStats[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $col = 2 + $j   # column B = 2
        $ws2.Cells.Item($row, $col).Value2 = $vals[$j]
    }
}

# Move the active selection on the "Knowledge Based Agent" sheet to D14.
$ws2.Activate()
[void]$ws2.Range("D14").Select()
